$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.253.72"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "2.283.39"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'113.71"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "'266.54"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").Value = "'47.77"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'9.35"
$ws.Range("E12").Value = "  +10.67%  "
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").Value = "'15.50"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "2.614.68"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "'0.868"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").Value = "2.288.67"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "43.311.71"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "'0.0000109"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "'6.86"
$ws.Range("E20").Value = "  +4.98%  "
$ws.Range("D21").Value = "'71.83"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("D23").Value = "'234.29"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("E24").Value = "  +3.77%  "
$ws.Range("D25").Value = "'2.91"
$ws.Range("E25").Value = "  +3.35%  "
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").Value = "'11.41"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "'4.00"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").Value = "'40.90"
$ws.Range("E29").Value = "  -3.77%  "
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").Value = "'173.71"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("D35").Value = "'5.74"
$ws.Range("E35").Value = "  +4.63%  "
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").Value = "'0.0370"
$ws.Range("E38").Value = "  +5.03%  "
$ws.Range("D39").Value = "'3.93"
$ws.Range("E39").Value = "  +3.70%  "
$ws.Range("E40").Value = "  -3.70%  "
$ws.Range("D41").Value = "'2.67"
$ws.Range("E41").Value = "  +10.52%  "
$ws.Range("D42").Value = "'77.81"
$ws.Range("E42").Value = "  +7.15%  "
$ws.Range("D43").Value = "'14.27"
$ws.Range("E43").Value = "  +6.14%  "
$ws.Range("D44").Value = "'0.239"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'6.28"
$ws.Range("E45").Value = "  +6.79%  "
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").Value = "'8.71"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").Value = "'104.70"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("E50").Value = "  +2.99%  "
$ws.Range("D51").Value = "'0.0997"
$ws.Range("E51").Value = "  -0.12%  "